$d = $word.ActiveDocument

$d.Content.Find.Execute("81÷9=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "48÷4=12, 0", 2) | Out-Null
$d.Content.Find.Execute("12÷6=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "56÷3=18, 2", 2) | Out-Null
$d.Content.Find.Execute("36÷9=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "88÷4=22, 0", 2) | Out-Null
$d.Content.Find.Execute("68÷3=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "65÷9=7, 2", 2) | Out-Null
$d.Content.Find.Execute("56÷9=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "27÷6=4, 3", 2) | Out-Null
$d.Content.Find.Execute("99÷6=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "86÷3=28, 2", 2) | Out-Null
$d.Content.Find.Execute("86÷6=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "49÷9=5, 4", 2) | Out-Null
$d.Content.Find.Execute("67÷3=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "29÷6=4, 5", 2) | Out-Null
$d.Content.Find.Execute("42÷4=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "94÷5=18, 4", 2) | Out-Null
$d.Content.Find.Execute("17÷8=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "30÷7=4, 2", 2) | Out-Null
$d.Content.Find.Execute("33÷2=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "93÷4=23, 1", 2) | Out-Null
$d.Content.Find.Execute("24÷6=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "92÷4=23, 0", 2) | Out-Null
$d.Content.Find.Execute("60÷6=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "72÷3=24, 0", 2) | Out-Null
$d.Content.Find.Execute("11÷4=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "97÷9=10, 7", 2) | Out-Null
$d.Content.Find.Execute("33÷4=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "37÷7=5, 2", 2) | Out-Null
$d.Content.Find.Execute("95÷5=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "12÷4=3, 0", 2) | Out-Null
$d.Content.Find.Execute("56÷6=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=7, 6", 2) | Out-Null
$d.Content.Find.Execute("18÷7=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "84÷7=12, 0", 2) | Out-Null
$d.Content.Find.Execute("42÷6=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "91÷9=10, 1", 2) | Out-Null
$d.Content.Find.Execute("43÷7=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "16÷3=5, 1", 2) | Out-Null
$d.Content.Find.Execute("39÷8=4, 7", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=4, 3", 2) | Out-Null
$d.Content.Find.Execute("45÷6=7, 3", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=6, 0", 2) | Out-Null
$d.Content.Find.Execute("44÷8=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=8, 4", 2) | Out-Null
$d.Content.Find.Execute("17÷5=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "39÷3=13, 0", 2) | Out-Null
$d.Content.Find.Execute("35÷4=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "58÷5=11, 3", 2) | Out-Null
